$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8072
$ws1.Range("F3").Value = 118
$ws1.Range("F4").Value = 92
$ws1.Range("F5").Value = 30520
$ws1.Range("F7").Value = 603
$ws1.Range("F8").Value = 694
$ws1.Range("F12").Value = 797
$ws1.Range("F13").Value = 51
$ws1.Range("F14").Value = 601
$ws1.Range("F15").Value = 381
$ws1.Range("F17").Value = 478
$ws1.Range("F19").Value = 411
$ws1.Range("F21").Value = 1112
$ws1.Range("F22").Value = 84
$ws1.Range("F23").Value = 698
$ws1.Range("F24").Value = 2327
$ws1.Range("F25").Value = 823
$ws1.Range("F26").Value = 65
$ws1.Range("F27").Value = 1085
$ws1.Range("F29").Value = 630
$ws1.Range("F30").Value = 1072

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 295
$ws2.Range("F4").Value = 339
$ws2.Range("F10").Value = 3

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 519

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 519
$ws4.Range("F3").Value = 8072
$ws4.Range("F4").Value = 118
$ws4.Range("F5").Value = 92
$ws4.Range("F6").Value = 295
$ws4.Range("F7").Value = 30521
$ws4.Range("F9").Value = 603
$ws4.Range("F10").Value = 694
$ws4.Range("F15").Value = 339
$ws4.Range("F18").Value = 797
$ws4.Range("F19").Value = 51
$ws4.Range("F20").Value = 601
$ws4.Range("F21").Value = 381
$ws4.Range("F26").Value = 3
$ws4.Range("F27").Value = 478
$ws4.Range("F29").Value = 411
$ws4.Range("F31").Value = 1112
$ws4.Range("F32").Value = 84
$ws4.Range("F33").Value = 698
$ws4.Range("F34").Value = 2327
$ws4.Range("F35").Value = 823
$ws4.Range("F36").Value = 65
$ws4.Range("F37").Value = 1085
$ws4.Range("F40").Value = 630
$ws4.Range("F41").Value = 1072

$wb.Save()
